# Daily attendance processing - 2026-01-28 14:03:37
# Swap the order of "Recorded By" (column G) entries from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colIndex = 7   # column G = "Recorded By"

$changed = 0
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $current = $cell.Value2
    if ($current -eq $oldText) {
        $cell.Value = $newText
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed 'Recorded By' cell(s) in column G."
